$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell updates per the cryptos-list refresh diff.
# Numeric-looking price strings get an explicit Text number format first so
# Excel's automatic type inference doesn't coerce them into floating-point numbers
# (which would lose formatting like trailing zeros or exact decimal text).

$ws.Range("D2").Value = "46.334.50"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "2.578.90"
$ws.Range("E3").Value = "  +9.42%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "304.51"
$ws.Range("E5").Value = "  +0.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "99.54"
$ws.Range("E6").Value = "  -0.25%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.599"
$ws.Range("E7").Value = "  +5.00%  "
$ws.Range("E9").Value = "  +12.68%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.51"
$ws.Range("E10").Value = "  +11.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0835"
$ws.Range("E11").Value = "  +4.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.13"
$ws.Range("E12").Value = "  +13.88%  "
$ws.Range("D13").Value = "2.973.06"
$ws.Range("E13").Value = "  +9.39%  "
$ws.Range("E14").Value = "  +1.85%  "
$ws.Range("D15").Value = "2.598.55"
$ws.Range("E15").Value = "  +10.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.900"
$ws.Range("E16").Value = "  +11.20%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.83"
$ws.Range("E17").Value = "  +8.70%  "
$ws.Range("D18").Value = "46.474.60"
$ws.Range("E18").Value = "  +0.79%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.32"
$ws.Range("E19").Value = "  +4.48%  "
$ws.Range("E20").Value = "  +4.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.64"
$ws.Range("E21").Value = "  +9.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.96"
$ws.Range("E22").Value = "  +5.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "254.92"
$ws.Range("E23").Value = "  +3.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.97"
$ws.Range("E24").Value = "  +4.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.17"
$ws.Range("E25").Value = "  +13.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "27.75"
$ws.Range("E26").Value = "  +32.04%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.17%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.45"
$ws.Range("E28").Value = "  +6.76%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.29"
$ws.Range("E29").Value = "  +4.61%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "39.56"
$ws.Range("E30").Value = "  -0.23%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.74"
$ws.Range("E31").Value = "  +0.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.11"
$ws.Range("E32").Value = "  +10.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.31"
$ws.Range("E33").Value = "  +22.85%  "
$ws.Range("E34").Value = "  +5.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0833"
$ws.Range("E35").Value = "  +7.43%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "149.64"
$ws.Range("E36").Value = "  +2.48%  "
$ws.Range("E37").Value = "  +3.70%  "
$ws.Range("E38").Value = "  +4.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.17"
$ws.Range("E39").Value = "  +5.43%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "15.75"
$ws.Range("E40").Value = "  +5.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.63"
$ws.Range("E41").Value = "  +12.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0322"
$ws.Range("E42").Value = "  +7.04%  "
$ws.Range("D43").Value = "2.028.64"
$ws.Range("E43").Value = "  +7.70%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "18.15"
$ws.Range("E45").Value = "  +25.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.38"
$ws.Range("E46").Value = "  -1.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.78"
$ws.Range("E47").Value = "  -0.77%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "108.48"
$ws.Range("E48").Value = "  +10.90%  "
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.04"
$ws.Range("E49").Value = "  +9.50%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.200"
$ws.Range("E50").Value = "  +7.46%  "
$ws.Range("D51").Value = "2.831.65"
$ws.Range("E51").Value = "  +9.51%  "
